$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.657.45"
$ws.Range("E2").Value = "  -2.34%  "
$ws.Range("D3").Value = "3.638.75"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.77%  "
$ws.Range("D7").Value = "3.632.15"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  -4.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +22.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.607"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "48.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("E14").Value = "  -3.04%  "
$ws.Range("D15").Value = "4.225.10"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "674.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "3.640.38"
$ws.Range("E18").Value = "  +2.03%  "
$ws.Range("D19").Value = "70.714.42"
$ws.Range("E19").Value = "  -2.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "99.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.59%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.90"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("E34").Value = "  -6.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "584.16"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.30"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.82%  "
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0455"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.14%  "
$ws.Range("D42").Value = "3.563.53"
$ws.Range("E42").Value = "  -3.24%  "
$ws.Range("E43").Value = "  -3.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.344"
$ws.Range("D44").Style = "Normal"
$ws.Range("E45").Value = "  -4.92%  "
$ws.Range("D46").Value = "0.0₃0734"
$ws.Range("E46").Value = "  -7.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.64%  "
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.03%  "
